$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cheltuieli si venituri")

# Update existing "Perioada" value for row 5 (Salariu) from "lunar" to "saptamanal"
$ws.Range("C5").Value = "saptamanal"

# Add a new row (row 8) for a weekly "Factura" (invoice) expense entry.
# First, copy formatting from the row above (A7/B7/C7/D7) so the new row
# keeps identical styling (date format, etc.) to the existing data rows.
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A8:D8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A8").Value = (Get-Date -Year 2012 -Month 4 -Day 11 -Hour 17 -Minute 3 -Second 10)
$ws.Range("B8").Value = "Factura"
$ws.Range("C8").Value = "saptamanal"
$ws.Range("D8").Value = 100

Write-Output "C5: $($ws.Range('C5').Value2)"
Write-Output "A8: $($ws.Range('A8').Value2) B8: $($ws.Range('B8').Value2) C8: $($ws.Range('C8').Value2) D8: $($ws.Range('D8').Value2)"
Write-Output "Dimension should now be A1:D8"
